$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2009")
$ws.Range("B2").Value = 0.5304054054054054
$ws.Range("B3").Value = 0.4182590233545648
$ws.Range("B4").Value = 0.2996870109546166
$ws.Range("B5").Value = 0.5895196506550219
$ws.Range("B6").Value = 0.397025171624714

$ws = $wb.Worksheets.Item("2018")
$ws.Range("B2").Value = 0.5470588235294118
$ws.Range("B3").Value = 0.4522184300341297
$ws.Range("B4").Value = 0.3742071881606766
$ws.Range("B5").Value = 0.5300261096605744
$ws.Range("B6").Value = 0.8142857142857143
$ws.Range("B7").Value = 0.5298245614035088
$ws.Range("B8").Value = 0.5176151761517616

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("C2").Value = 0.4469792523988644
$ws.Range("D2").Value = 0.1143382637970692
$ws.Range("E2").Value = 0.2996870109546166
$ws.Range("F2").Value = 0.397025171624714
$ws.Range("G2").Value = 0.4182590233545648
$ws.Range("H2").Value = 0.5304054054054054
$ws.Range("I2").Value = 0.5895196506550219
$ws.Range("C3").Value = 0.566727398310069
$ws.Range("D3").Value = 0.07901800925838894
$ws.Range("E3").Value = 0.4843373493975904
$ws.Range("F3").Value = 0.5014942663006445
$ws.Range("G3").Value = 0.5532683279162152
$ws.Range("H3").Value = 0.6047430830039525
$ws.Range("I3").Value = 0.7132352941176471
$ws.Range("C4").Value = 0.4919070526162364
$ws.Range("D4").Value = 0.08588513811802025
$ws.Range("E4").Value = 0.353874883286648
$ws.Range("F4").Value = 0.4433497536945813
$ws.Range("G4").Value = 0.5014245014245015
$ws.Range("H4").Value = 0.5204819277108433
$ws.Range("I4").Value = 0.6305732484076433
$ws.Range("C5").Value = 0.4716381970218835
$ws.Range("D5").Value = 0.08730082396682598
$ws.Range("E5").Value = 0.3126069594980034
$ws.Range("F5").Value = 0.4186915887850468
$ws.Range("G5").Value = 0.4684642686867617
$ws.Range("H5").Value = 0.5185674009784682
$ws.Range("I5").Value = 0.631578947368421
$ws.Range("C6").Value = 0.5090185306926361
$ws.Range("D6").Value = 0.09024197937760119
$ws.Range("E6").Value = 0.3520286396181384
$ws.Range("F6").Value = 0.4417177914110429
$ws.Range("G6").Value = 0.5072886297376094
$ws.Range("H6").Value = 0.5571428571428572
$ws.Range("I6").Value = 0.7074829931972789
$ws.Range("C7").Value = 0.4877116648806603
$ws.Range("D7").Value = 0.1228332128482508
$ws.Range("E7").Value = 0.3367088607594937
$ws.Range("F7").Value = 0.372302312407402
$ws.Range("G7").Value = 0.5064102564102564
$ws.Range("H7").Value = 0.5842172934747669
$ws.Range("I7").Value = 0.6571428571428571
$ws.Range("C8").Value = 0.4850568301597232
$ws.Range("D8").Value = 0.07179638267937676
$ws.Range("E8").Value = 0.3932432432432432
$ws.Range("F8").Value = 0.4240440248836299
$ws.Range("G8").Value = 0.4976851851851852
$ws.Range("H8").Value = 0.5419003438805419
$ws.Range("I8").Value = 0.5725806451612904
$ws.Range("C9").Value = 0.5255727463292839
$ws.Range("D9").Value = 0.08588126891981591
$ws.Range("E9").Value = 0.3980099502487562
$ws.Range("F9").Value = 0.4463362068965517
$ws.Range("G9").Value = 0.5354654441594353
$ws.Range("H9").Value = 0.5800148783476282
$ws.Range("I9").Value = 0.6692913385826772
$ws.Range("C10").Value = 0.46997086920851
$ws.Range("D10").Value = 0.1140574647312718
$ws.Range("E10").Value = 0.2226613965744401
$ws.Range("F10").Value = 0.4006610576923077
$ws.Range("G10").Value = 0.4705067897752104
$ws.Range("H10").Value = 0.5469582629993712
$ws.Range("I10").Value = 0.6687898089171974
$ws.Range("C11").Value = 0.5378908576036824
$ws.Range("D11").Value = 0.136105534751788
$ws.Range("E11").Value = 0.3742071881606766
$ws.Range("F11").Value = 0.4849168030929456
$ws.Range("G11").Value = 0.5298245614035088
$ws.Range("H11").Value = 0.5385424665949932
$ws.Range("I11").Value = 0.8142857142857143

$ws = $wb.Worksheets.Item("2010")
$ws.Range("B2").Value = 0.5942028985507246
$ws.Range("B3").Value = 0.568075117370892
$ws.Range("B4").Value = 0.5034168564920274
$ws.Range("B5").Value = 0.4843373493975904
$ws.Range("B6").Value = 0.5384615384615384
$ws.Range("B7").Value = 0.7132352941176471
$ws.Range("B8").Value = 0.4957264957264957
$ws.Range("B9").Value = 0.6363636363636364

$ws = $wb.Worksheets.Item("2011")
$ws.Range("B2").Value = 0.5136612021857924
$ws.Range("B3").Value = 0.4838709677419355
$ws.Range("B4").Value = 0.5204819277108433
$ws.Range("B5").Value = 0.5014245014245015
$ws.Range("B6").Value = 0.5859872611464968
$ws.Range("B7").Value = 0.356475300400534
$ws.Range("B8").Value = 0.353874883286648
$ws.Range("B9").Value = 0.4989517819706499
$ws.Range("B10").Value = 0.3695090439276486
$ws.Range("B11").Value = 0.5157894736842106
$ws.Range("B12").Value = 0.4433497536945813
$ws.Range("B13").Value = 0.4668769716088328
$ws.Range("B14").Value = 0.6305732484076433
$ws.Range("B15").Value = 0.4040747028862479
$ws.Range("B16").Value = 0.6083333333333333
$ws.Range("B17").Value = 0.5134408602150538
$ws.Range("B18").Value = 0.5957446808510638

$ws = $wb.Worksheets.Item("2012")
$ws.Range("B2").Value = 0.5618729096989966
$ws.Range("B3").Value = 0.4205607476635514
$ws.Range("B4").Value = 0.4510556621880998
$ws.Range("B5").Value = 0.3126069594980034
$ws.Range("B6").Value = 0.5041322314049587
$ws.Range("B7").Value = 0.5632183908045977
$ws.Range("B8").Value = 0.4578587699316629
$ws.Range("B9").Value = 0.4851258581235698
$ws.Range("B10").Value = 0.3794940079893475
$ws.Range("B11").Value = 0.631578947368421
$ws.Range("B12").Value = 0.4130841121495327
$ws.Range("B13").Value = 0.4790697674418605

$ws = $wb.Worksheets.Item("2013")
$ws.Range("B2").Value = 0.5171102661596958
$ws.Range("B3").Value = 0.4111111111111111
$ws.Range("B4").Value = 0.3520286396181384
$ws.Range("B5").Value = 0.6641221374045801
$ws.Range("B6").Value = 0.5533596837944664
$ws.Range("B7").Value = 0.5718390804597702
$ws.Range("B8").Value = 0.4651162790697674
$ws.Range("B9").Value = 0.7074829931972789
$ws.Range("B10").Value = 0.4680851063829787
$ws.Range("B11").Value = 0.4376321353065539
$ws.Range("B12").Value = 0.4417177914110429
$ws.Range("B13").Value = 0.5048076923076923
$ws.Range("B14").Value = 0.5571428571428572
$ws.Range("B15").Value = 0.5072886297376094
$ws.Range("B16").Value = 0.5110132158590308
$ws.Range("B17").Value = 0.4193548387096774
$ws.Range("B18").Value = 0.5641025641025641

$ws = $wb.Worksheets.Item("2014")
$ws.Range("B2").Value = 0.5064102564102564
$ws.Range("B3").Value = 0.3646551724137931
$ws.Range("B4").Value = 0.5409836065573771
$ws.Range("B5").Value = 0.5376344086021505
$ws.Range("B6").Value = 0.4030837004405287
$ws.Range("B7").Value = 0.6274509803921569
$ws.Range("B8").Value = 0.6461538461538462
$ws.Range("B9").Value = 0.6571428571428571
$ws.Range("B10").Value = 0.379949452401011
$ws.Range("B11").Value = 0.3367088607594937
$ws.Range("B12").Value = 0.3646551724137931

$ws = $wb.Worksheets.Item("2015")
$ws.Range("B2").Value = 0.3932432432432432
$ws.Range("B3").Value = 0.5313531353135313
$ws.Range("B4").Value = 0.5725806451612904
$ws.Range("B5").Value = 0.4074803149606299
$ws.Range("B6").Value = 0.4406077348066298
$ws.Range("B7").Value = 0.5524475524475524
$ws.Range("B8").Value = 0.4976851851851852

$ws = $wb.Worksheets.Item("2016")
$ws.Range("B2").Value = 0.631336405529954
$ws.Range("B3").Value = 0.5609756097560976
$ws.Range("B4").Value = 0.5786802030456852
$ws.Range("B5").Value = 0.6136363636363636
$ws.Range("B6").Value = 0.4367816091954023
$ws.Range("B7").Value = 0.4983277591973244
$ws.Range("B8").Value = 0.545816733067729
$ws.Range("B9").Value = 0.3980099502487562
$ws.Range("B10").Value = 0.475
$ws.Range("B11").Value = 0.4215116279069768
$ws.Range("B12").Value = 0.5251141552511416
$ws.Range("B13").Value = 0.6692913385826772
$ws.Range("B14").Value = 0.5804597701149425
$ws.Range("B15").Value = 0.4230769230769231

$ws = $wb.Worksheets.Item("2017")
$ws.Range("B2").Value = 0.4211886304909561
$ws.Range("B3").Value = 0.6063348416289592
$ws.Range("B4").Value = 0.4385499557913351
$ws.Range("B5").Value = 0.6687898089171974
$ws.Range("B6").Value = 0.4893048128342246
$ws.Range("B7").Value = 0.4026442307692308
$ws.Range("B8").Value = 0.5268817204301075
$ws.Range("B9").Value = 0.3858064516129032
$ws.Range("B10").Value = 0.3645484949832776
$ws.Range("B11").Value = 0.2226613965744401
$ws.Range("B12").Value = 0.551440329218107
$ws.Range("B13").Value = 0.5015384615384615
$ws.Range("B14").Value = 0.4
$ws.Range("B15").Value = 0.588
$ws.Range("B16").Value = 0.4517087667161961
$ws.Range("B17").Value = 0.3105065666041276
$ws.Range("B18").Value = 0.5960591133004927
$ws.Range("B19").Value = 0.5335120643431636
